$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1062
$ws.Range("F3").Value = 674
$ws.Range("F4").Value = 1489
$ws.Range("F5").Value = 3246
$ws.Range("F7").Value = 666
$ws.Range("F8").Value = 2221
$ws.Range("F9").Value = 479
$ws.Range("F11").Value = 235
$ws.Range("F12").Value = 127
$ws.Range("F13").Value = 310
$ws.Range("F14").Value = 1068
$ws.Range("F15").Value = 430
$ws.Range("F16").Value = 6
$ws.Range("F18").Value = 206
$ws.Range("F19").Value = 4452
$ws.Range("F20").Value = 1294
$ws.Range("F21").Value = 3378
$ws.Range("F23").Value = 75
$ws.Range("F24").Value = 165
$ws.Range("F25").Value = 3314
$ws.Range("F26").Value = 4928
$ws.Range("F29").Value = 543
$ws.Range("F30").Value = 3186
$ws.Range("F31").Value = 349
$ws.Range("F35").Value = 872
$ws.Range("F36").Value = 1157
$ws.Range("F37").Value = 1398
$ws.Range("F38").Value = 114
$ws.Range("F39").Value = 1325
$ws.Range("F40").Value = 843
$ws.Range("F42").Value = 796
$ws.Range("F45").Value = 287
$ws.Range("F46").Value = 59
$ws.Range("F47").Value = 144
$ws.Range("F49").Value = 3710

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 997
$ws.Range("F11").Value = 7
$ws.Range("F21").Value = 39

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2107

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2107
$ws.Range("F3").Value = 674
$ws.Range("F4").Value = 1489
$ws.Range("F5").Value = 3246
$ws.Range("F7").Value = 666
$ws.Range("F9").Value = 2221
$ws.Range("F10").Value = 479
$ws.Range("F12").Value = 235
$ws.Range("F13").Value = 997
$ws.Range("F14").Value = 127
$ws.Range("F15").Value = 310
$ws.Range("F16").Value = 1068
$ws.Range("F17").Value = 430
$ws.Range("F18").Value = 6
$ws.Range("F19").Value = 206
$ws.Range("F20").Value = 4452
$ws.Range("F21").Value = 1294
$ws.Range("F23").Value = 3378
$ws.Range("F24").Value = 3314
$ws.Range("F25").Value = 4928
$ws.Range("F28").Value = 3186
$ws.Range("F29").Value = 349
$ws.Range("F33").Value = 872
$ws.Range("F34").Value = 1157
$ws.Range("F35").Value = 1398
$ws.Range("F36").Value = 114
$ws.Range("F37").Value = 1325
$ws.Range("F39").Value = 844
$ws.Range("F44").Value = 287
$ws.Range("F46").Value = 59
$ws.Range("F47").Value = 144
$ws.Range("F48").Value = 363
$ws.Range("F49").Value = 3711
